# Append new clock in/out log entries to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$entries = @(
    "OUT -> 2017/03/23 15:49",
    "IN -> 2017/03/23 15:50",
    "OUT -> 2017/03/23 15:52",
    "IN -> 2017/03/23 15:53",
    "OUT -> 2017/03/23 15:53",
    "IN -> 2017/03/23 15:55",
    "OUT -> 2017/03/27 14:55",
    "IN -> 2017/03/27 14:57",
    "OUT -> 2017/03/27 15:02",
    "IN -> 2017/03/27 15:02",
    "OUT -> 2017/03/27 15:04",
    "IN -> 2017/03/27 15:05",
    "OUT -> 2017/03/27 15:07",
    "IN -> 2017/03/27 15:08",
    "OUT -> 2017/03/27 15:11",
    "IN -> 2017/03/27 15:11",
    "OUT -> 2017/03/27 15:14",
    "IN -> 2017/03/27 15:15",
    "OUT -> 2017/03/27 15:16",
    "IN -> 2017/03/27 15:17",
    "OUT -> 2017/03/27 15:22",
    "IN -> 2017/03/27 15:25"
)

# Continue the log right after the last populated row in column A.
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
$startRow = $lastRow + 1

for ($i = 0; $i -lt $entries.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $entries[$i]
}
